# Updated cryptos list refresh (prices / 1h volume change %) plus a
# 3-way reshuffle of the Hedera / InjectiveProtocol / OKB rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "64.657.99"
$ws.Range("E2").Value = "  -2.34%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "3.138.23"
$ws.Range("E3").Value = "  -8.09%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.06%  "

# --- Row 5: BNB ---
$ws.Range("D5").Value = "564.08"
$ws.Range("E5").Value = "  -3.30%  "

# --- Row 6: Solana ---
$ws.Range("D6").Value = "169.95"
$ws.Range("E6").Value = "  -4.71%  "

# --- Row 7: XRP ---
$ws.Range("D7").Value = "0.615"
$ws.Range("E7").Value = "  -1.17%  "

# --- Row 8: USDC ---
$ws.Range("E8").Value = "  +0.04%  "

# --- Row 9: LidoStakedEther ---
$ws.Range("D9").Value = "3.133.43"
$ws.Range("E9").Value = "  -8.23%  "

# --- Row 10: Dogecoin ---
$ws.Range("E10").Value = "  -5.95%  "

# --- Row 11: Toncoin ---
$ws.Range("D11").Value = "6.55"
$ws.Range("E11").Value = "  -5.89%  "

# --- Row 12: Cardano ---
$ws.Range("D12").Value = "0.392"
$ws.Range("E12").Value = "  -4.97%  "

# --- Row 13: WrappedliquidstakedEther2.0 ---
$ws.Range("D13").Value = "3.677.89"
$ws.Range("E13").Value = "  -8.24%  "

# --- Row 14: TRON ---
$ws.Range("E14").Value = "  +0.91%  "

# --- Row 15: Avalanche ---
$ws.Range("D15").Value = "26.97"
$ws.Range("E15").Value = "  -7.99%  "

# --- Row 16: WrappedBTC ---
$ws.Range("D16").Value = "64.538.47"
$ws.Range("E16").Value = "  -2.71%  "

# --- Row 17: ShibaInu ---
$ws.Range("D17").Value = "0.0000162"
$ws.Range("E17").Value = "  -5.85%  "

# --- Row 18: WrappedEther ---
$ws.Range("D18").Value = "3.140.93"

# --- Row 19: Polkadot ---
$ws.Range("D19").Value = "5.67"
$ws.Range("E19").Value = "  -3.88%  "

# --- Row 20: Chainlink ---
$ws.Range("D20").Value = "12.83"
$ws.Range("E20").Value = "  -6.72%  "

# --- Row 21: BitcoinCash ---
$ws.Range("D21").Value = "354.58"
$ws.Range("E21").Value = "  -3.17%  "

# --- Row 22: Uniswap ---
$ws.Range("E22").Value = "  -4.71%  "

# --- Row 23: Dai ---
$ws.Range("E23").Value = "  +0.47%  "

# --- Row 24: Litecoin ---
$ws.Range("D24").Value = "68.35"
$ws.Range("E24").Value = "  -6.03%  "

# --- Row 25: PEPE ---
$ws.Range("D25").Value = "0.0000117"
$ws.Range("E25").Value = "  -6.75%  "

# --- Row 26: Polygon ---
$ws.Range("D26").Value = "0.498"
$ws.Range("E26").Value = "  -6.87%  "

# --- Row 27: InternetComputer(DFINITY) ---
$ws.Range("D27").Value = "9.56"
$ws.Range("E27").Value = "  -2.56%  "

# --- Row 28: Kaspa ---
$ws.Range("E28").Value = "  -2.60%  "

# --- Row 29: Binance-PegBSC-USD ---
# Force text so Excel keeps the trailing zero ("1.00", not "1").
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.12%  "

# --- Row 30: USDe ---
$ws.Range("E30").Value = "  -0.16%  "

# --- Row 31: PancakeSwap ---
$ws.Range("E31").Value = "  -4.94%  "

# --- Row 32: NEARProtocol ---
$ws.Range("D32").Value = "5.35"
$ws.Range("E32").Value = "  -7.18%  "

# --- Row 33: EthereumClassic ---
# Force text so Excel keeps the trailing zero ("21.90", not "21.9").
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "21.90"
$ws.Range("E33").Value = "  -6.40%  "

# --- Row 34: Aptos ---
$ws.Range("D34").Value = "6.62"
$ws.Range("E34").Value = "  -5.59%  "

# --- Row 35: Fetch.AI ---
$ws.Range("D35").Value = "1.19"
$ws.Range("E35").Value = "  -5.19%  "

# --- Row 36: ImmutableX ---
$ws.Range("D36").Value = "1.43"
$ws.Range("E36").Value = "  -7.22%  "

# --- Row 37: Monero ---
$ws.Range("D37").Value = "153.11"
$ws.Range("E37").Value = "  -5.82%  "

# --- Row 38: Mantle ---
$ws.Range("D38").Value = "0.828"
$ws.Range("E38").Value = "  -5.25%  "

# --- Row 39: EnergySwap ---
$ws.Range("D39").Value = "25.99"
$ws.Range("E39").Value = "  -5.89%  "

# --- Row 40: Stacks ---
$ws.Range("D40").Value = "1.73"
$ws.Range("E40").Value = "  -3.00%  "

# --- Row 41: dogwifhat ---
$ws.Range("D41").Value = "2.52"
$ws.Range("E41").Value = "  -1.93%  "

# --- Row 42: Maker ---
$ws.Range("D42").Value = "2.648.77"
$ws.Range("E42").Value = "  -1.81%  "

# --- Row 43: Filecoin ---
$ws.Range("D43").Value = "4.16"
$ws.Range("E43").Value = "  -6.38%  "

# --- Row 44: RenderToken ---
$ws.Range("D44").Value = "6.01"
$ws.Range("E44").Value = "  -5.56%  "

# --- Rows 45-47: three-way reshuffle.
# Old order: 45=Hedera, 46=InjectiveProtocol, 47=OKB
# New order: 45=InjectiveProtocol, 46=OKB, 47=Hedera
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "24.09"
$ws.Range("E45").Value = "  -4.48%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
# Force text so Excel keeps the trailing zero ("39.00", not "39").
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.00"
$ws.Range("E46").Value = "  -2.27%  "

$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "0.0652"
$ws.Range("E47").Value = "  -5.21%  "

# --- Row 48: Bittensor ---
$ws.Range("D48").Value = "319.09"
$ws.Range("E48").Value = "  -4.75%  "

# --- Row 49: VeChain ---
$ws.Range("D49").Value = "0.0272"
$ws.Range("E49").Value = "  -4.65%  "

# --- Row 50: Stellar ---
$ws.Range("E50").Value = "  -2.36%  "

# --- Row 51: FirstDigitalUSD ---
$ws.Range("E51").Value = "  -0.14%  "
